$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A11').Value = 'آلاخان والاخان'
$ws.Range('B11').Value = 'آلاخان والاخاني جأ بئتره'
$ws.Range('C11').Value = 'دربه‌در'
$ws.Range('D11').Value = 'از دربه‌دری بهتر است'
$ws.Range('A12').Value = 'آها'
$ws.Range('B12').Value = 'آها، منم أيم'
$ws.Range('C12').Value = 'بله'
$ws.Range('D12').Value = 'بله، من هم می‌آیم'
$ws.Range('A13').Value = 'آهين'
$ws.Range('B13').Value = 'مي دۊکان آهين پۊرد ٚ ور نهأ'
$ws.Range('C13').Value = 'آهن'
$ws.Range('D13').Value = 'مغازه من کنار پل آهنی است'
$ws.Range('A14').Value = 'آوج'
$ws.Range('B14').Value = 'هرچي زنگ بزئم آوجأ ندأ'
$ws.Range('C14').Value = 'جواب'
$ws.Range('D14').Value = 'هر چه زنگ زدم جواب نداد'
$ws.Range('A15').Value = 'ارجئن ارجئن'
$ws.Range('B15').Value = '﻿کاغذانأ همهˈ ارجئن ارجئنأ کۊده'
$ws.Range('C15').Value = 'ریز ریز'
$ws.Range('D15').Value = 'تمام کاغذها را ریز ریز کرد'
$ws.Range('A16').Value = 'اي'
$ws.Range('B16').Value = 'اي ميليۊن تۊمۊن پۊل خأیم'
$ws.Range('C16').Value = 'یک'
$ws.Range('D16').Value = 'یک میلیون تومان پول لازم دارم'
$ws.Range('A17').Value = 'اي فچم'
$ws.Range('B17').Value = 'ترأ بيدينم اي فچم ماچي دهم'
$ws.Range('C17').Value = 'بسیار زیاد'
$ws.Range('D17').Value = 'تو را ببینم بسیار زیاد مى‌بوسم‌ات'
$ws.Range('A18').Value = 'اي وخت'
$ws.Range('B18').Value = 'اي وخت پۊر ترأ خأستيم'
$ws.Range('C18').Value = 'زمانی'
$ws.Range('D18').Value = 'زمانی تو را خیلى می‌خواستم'
$ws.Range('A19').Value = 'اي ور'
$ws.Range('B19').Value = 'أن ٚ رختان هرتأ اي ور کفته'
$ws.Range('C19').Value = 'یک گوشه'
$ws.Range('D19').Value = 'هرکدام از لباس‌هایش یک گوشه افتاده است'
$ws.Range('A20').Value = 'اي وري'
$ws.Range('B20').Value = 'أ تابلؤ چي ره اي وري نهأ؟'
$ws.Range('C20').Value = 'کج'
$ws.Range('D20').Value = 'برای چه این تابلو کج است؟'
$ws.Range('A21').Value = 'اي پئم'
$ws.Range('B21').Value = 'اي پئم سيبيشکا مرأ فأدي؟'
$ws.Range('C21').Value = 'یک مشت'
$ws.Range('D21').Value = 'یک مشت تخمه آفتابگردان به من مى‌دهى؟'
$ws.Range('A22').Value = 'ايتأ'
$ws.Range('B22').Value = 'ايتأ کاغذ مرأ فأدن'
$ws.Range('C22').Value = 'یک'
$ws.Range('D22').Value = 'یک کاغذ به من بده'
$ws.Range('A23').Value = 'ايجانا'
$ws.Range('B23').Value = 'من ؤ ميعاد ايجانا بۊشؤئيم بيرۊن'
$ws.Range('C23').Value = 'با هم'
$ws.Range('D23').Value = 'من و میعاد با هم بیرون رفتیم'
$ws.Range('A24').Value = 'ايدانه'
$ws.Range('B24').Value = 'ايدانه کۊئي خأيم'
$ws.Range('C24').Value = 'یک عدد'
$ws.Range('D24').Value = 'یک عدد کدو مى‌خواهم'
$ws.Range('A25').Value = 'ايدفأ'
$ws.Range('B25').Value = 'ايدفأ بامؤم شيمي خانه'
$ws.Range('C25').Value = 'یک بار'
$ws.Range('D25').Value = 'یک بار به خانه شما آمدم'
$ws.Range('A26').Value = 'ايدفأئي'
$ws.Range('B26').Value = 'ناصر ايدفأئي بامؤ اۊتاق ٚ دۊرۊن'
$ws.Range('C26').Value = 'ناگهان'
$ws.Range('D26').Value = 'ناصر ناگهان داخل اتاق آمد'
$ws.Range('A27').Value = 'ايدقه'
$ws.Range('B27').Value = 'ايدقه بئس'
$ws.Range('C27').Value = 'یک دقیقه'
$ws.Range('D27').Value = 'یک دقیقه صبر کن'
$ws.Range('A28').Value = 'ايرۊز'
$ws.Range('B28').Value = 'ايرۊز دئه بمانسته مي تعطيلي-يان تۊمانأ به'
$ws.Range('C28').Value = 'یک روز'
$ws.Range('D28').Value = 'یک روز دیگر مانده تا تعطیلاتم تمام شود'
$ws.Range('A29').Value = 'ايزه'
$ws.Range('B29').Value = 'ايزه بۊشۊ اۊشنتر'
$ws.Range('C29').Value = 'کمى'
$ws.Range('D29').Value = 'کمى آن طرف‌تر برو'
$ws.Range('A30').Value = 'ايسأبيد'
$ws.Range('B30').Value = 'ديشب همه أمي خانه ايسابيد'
$ws.Range('C30').Value = 'بودند'
$ws.Range('D30').Value = 'همه دیشب خانه ما بودند'
$ws.Range('A31').Value = 'ايستۊل'
$ws.Range('B31').Value = 'أشأن چۊب ٚ جأ ایستۊل چأکۊنيد'
$ws.Range('C31').Value = 'صندلی'
$ws.Range('D31').Value = 'آن‌ها از چوب صندلی می‌سازند'
$ws.Range('A32').Value = 'ايسم'
$ws.Range('B32').Value = 'تي ايسم چيسه؟'
$ws.Range('C32').Value = 'نام'
$ws.Range('D32').Value = 'نام تو چیست؟'
$ws.Range('A33').Value = 'ايسيد؟'
$ws.Range('B33').Value = 'شۊمان چن نفر ايسيد؟'
$ws.Range('C33').Value = 'هستید؟'
$ws.Range('D33').Value = 'شما چند نفر هستید؟'
$ws.Range('A34').Value = 'ايسکمۊ'
$ws.Range('B34').Value = 'بۊشۊ دۊکان دۊتا ايسکمۊ بيهين'
$ws.Range('C34').Value = 'بستنی یخی'
$ws.Range('D34').Value = 'برو مغازه دو عدد بستنی یخی بخر'
$ws.Range('A35').Value = 'ايشتاو'
$ws.Range('B35').Value = 'هأی أمرأ ايشتاو ديهه'
$ws.Range('C35').Value = 'هشدار'
$ws.Range('D35').Value = 'مدام به ما هشدار می‌دهد'
$ws.Range('A36').Value = 'ايشتاوي'
$ws.Range('B36').Value = 'ايشتاوي چي گمه؟'
$ws.Range('C36').Value = 'می‌شنوی'
$ws.Range('D36').Value = 'می‌شنوی چه می‌گویم؟'
$ws.Range('A37').Value = 'ايششه'
$ws.Range('B37').Value = 'ايششه، لانتي'
$ws.Range('C37').Value = 'اه (لفظ تنفر)'
$ws.Range('D37').Value = 'اه، مار'
$ws.Range('A38').Value = 'ايشماردن'
$ws.Range('B38').Value = 'مي زأى ايشماردنأ نأنه'
$ws.Range('C38').Value = 'شمارش'
$ws.Range('D38').Value = 'بچه من شمارش بلد نیست'
$ws.Range('A39').Value = 'ايشکؤر'
$ws.Range('B39').Value = 'کيشکائان ٚ رئم ايشکؤر فۊکۊن'
$ws.Range('C39').Value = 'خرده‌برنج'
$ws.Range('D39').Value = 'برای جوجه‌ها هم خرده‌برنج بریز'
$ws.Range('A40').Value = 'ايشکنئن'
$ws.Range('B40').Value = 'چي-يأ ايشکنئن دري؟'
$ws.Range('C40').Value = 'شکستن'
$ws.Range('D40').Value = 'دارى چه چیزى را مى‌شکنى؟'
$ws.Range('A41').Value = 'ايمرۊ'
$ws.Range('B41').Value = 'ايمرۊ بۊشؤم کيتابخانه'
$ws.Range('C41').Value = 'امروز'
$ws.Range('D41').Value = 'امروز به کتابخانه رفتم'
$ws.Range('A42').Value = 'اينسان'
$ws.Range('B42').Value = 'اينسان بۊبۊ'
$ws.Range('C42').Value = 'انسان'
$ws.Range('D42').Value = 'انسان باش'
$ws.Range('A43').Value = 'ايوار ايوار'
$ws.Range('B43').Value = 'ايوار ايوار شئطان مرأ گه…'
$ws.Range('C43').Value = 'گهگاه'
$ws.Range('D43').Value = 'گهگاه شیطان به من مى‌گوید…'
$ws.Range('A44').Value = 'ايواردم'
$ws.Range('B44').Value = 'ايواردم وا بشم بازار'
$ws.Range('C44').Value = 'یک بار دیگر'
$ws.Range('D44').Value = 'یک بار دیگر باید به بازار بروم'
$ws.Range('A45').Value = 'ايوارٚکي'
$ws.Range('B45').Value = 'ايوارٚکي بۊگۊ نأيم ترأ راحتأ کۊن'
$ws.Range('C45').Value = 'به یک باره'
$ws.Range('D45').Value = 'به یک باره بگو نمى‌آیم خودت را خلاص کن'
$ws.Range('A46').Value = 'ايپيچم'
$ws.Range('B46').Value = 'گۊش نۊکۊني ايپيچم ترأيأ'
$ws.Range('C46').Value = 'برخورد تندی می‌کنم'
$ws.Range('D46').Value = 'توجه نکنى برخورد تندى با تو مى‌کنم'
$ws.Range('A47').Value = 'ايپچه'
$ws.Range('B47').Value = 'ايپچه پلا مه ره دۊکۊن'
$ws.Range('C47').Value = 'کمی'
$ws.Range('D47').Value = 'کمى پلو براى من بریز'
$ws.Range('A48').Value = 'ايچي'
$ws.Range('B48').Value = 'ايچي خأيم ترأ بگم'
$ws.Range('C48').Value = 'یک‌چیز'
$ws.Range('D48').Value = 'یک‌چیز مى‌خواهم به تو بگویم'
$ws.Range('A49').Value = 'اۊتاق'
$ws.Range('B49').Value = 'مي اۊتاق کؤيتأ ايسه؟'
$ws.Range('C49').Value = 'اتاق'
$ws.Range('D49').Value = 'اتاق من کدام یک است؟'
$ws.Range('A50').Value = 'اۊدۊشتن'
$ws.Range('B50').Value = 'سۊمبۊر ٚ مأنستن أمي خۊنأ اۊدۊشتن دريد'
$ws.Range('C50').Value = 'مکیدن'
$ws.Range('D50').Value = 'مانند زالو دارند خون ما را مى‌دوشند'
$ws.Range('A51').Value = 'اۊرشين'
$ws.Range('B51').Value = 'کترأيأ اۊسادم قاتؤقا اۊرشين بزئم'
$ws.Range('C51').Value = 'هم زدن'
$ws.Range('D51').Value = 'کفگیر چوبى را برداشتم خورش را هم زدم'
$ws.Range('A52').Value = 'اۊسان'
$ws.Range('B52').Value = 'قاقؤشاقانأ اۊسان بأور'
$ws.Range('C52').Value = 'بردار'
$ws.Range('D52').Value = 'قاشق‌ها را بردار بیاور'
$ws.Range('A53').Value = 'اۊسه کۊد'
$ws.Range('B53').Value = 'أ کيتابأ کي اۊسه کۊد؟'
$ws.Range('C53').Value = 'فرستاد'
$ws.Range('D53').Value = 'این کتاب را چه کسی فرستاد؟'
$ws.Range('A54').Value = 'اۊشان'
$ws.Range('B54').Value = 'اۊشان کيسيد؟'
$ws.Range('C54').Value = 'آن‌ها'
$ws.Range('D54').Value = 'آن‌ها کیستند؟'
$ws.Range('A55').Value = 'اۊشکؤفت'
$ws.Range('B55').Value = 'آب مرأ اۊشکؤفت'
$ws.Range('C55').Value = 'در گلو گیر کرد'
$ws.Range('D55').Value = 'آب در گلویم گیر کرد'
$ws.Range('A56').Value = 'اۊن'
$ws.Range('B56').Value = 'اۊن چيسه؟'
$ws.Range('C56').Value = 'آن'
$ws.Range('D56').Value = 'آن چیست؟'
$ws.Range('A57').Value = 'اۊيأ'
$ws.Range('B57').Value = 'مهين ٚ کيف اۊيأ نهأ'
$ws.Range('C57').Value = 'آن‌جا'
$ws.Range('D57').Value = 'کیف مهین آن‌جا است'
$ws.Range('A58').Value = 'اۊچين'
$ws.Range('B58').Value = 'تي کيتابانأ اۊچين'
$ws.Range('C58').Value = 'جمع کن'
$ws.Range('D58').Value = 'کتاب‌هایت را جمع کن'
$ws.Range('A59').Value = 'اۊکۊف'
$ws.Range('B59').Value = 'أمه ره اۊکۊف دره'
$ws.Range('C59').Value = 'شگون'
$ws.Range('D59').Value = 'براى ما شگون دارد'
$ws.Range('A60').Value = 'أ'
$ws.Range('B60').Value = 'شۊفاژأ چأکۊدي؟'
$ws.Range('C60').Value = 'را'
$ws.Range('D60').Value = 'شوفاژ را تعمیر کردی؟'
$ws.Range('A61').Value = 'أتؤ'
$ws.Range('B61').Value = 'أتؤ نيه'
$ws.Range('C61').Value = 'چنین'
$ws.Range('D61').Value = 'چنین نیست'
$ws.Range('A62').Value = 'أجۊر'
$ws.Range('B62').Value = 'أجۊر چي-يان مرأ حالي نيه'
$ws.Range('C62').Value = 'این‌گونه'
$ws.Range('D62').Value = 'این‌گونه چیزها را متوجه نمى‌شوم'
$ws.Range('A63').Value = 'أذب'
$ws.Range('B63').Value = 'مي برار أذبه'
$ws.Range('C63').Value = 'مجرد'
$ws.Range('D63').Value = 'برادرم مجرد است'
$ws.Range('A64').Value = 'أرايم'
$ws.Range('B64').Value = 'أرايم بيأ'
$ws.Range('C64').Value = 'این‌طرف هم'
$ws.Range('D64').Value = 'این‌طرف هم بیا'
$ws.Range('A65').Value = 'أرسۊ'
$ws.Range('B65').Value = 'تي أرسۊئانأ پاکۊن'
$ws.Range('C65').Value = 'اشک'
$ws.Range('D65').Value = 'اشک‌‌هایت را پاک کن'
$ws.Range('A66').Value = 'أزازيل'
$ws.Range('B66').Value = 'أ أزازيل ٚ زاکانأ فأندر'
$ws.Range('C66').Value = 'بسیار شیطان'
$ws.Range('D66').Value = 'این بچه‌هاى بسیار شیطان را نگاه کن'
$ws.Range('A67').Value = 'أسباب'
$ws.Range('B67').Value = 'ايپچه أسباب بأورم بۊخۊري؟'
$ws.Range('C67').Value = 'خوراکى'
$ws.Range('D67').Value = 'کمى خوراکى بیاورم بخورى؟'
$ws.Range('A68').Value = 'أشان'
$ws.Range('B68').Value = 'أشان کي أمؤن دريد کيسيد؟'
$ws.Range('C68').Value = 'ایشان'
$ws.Range('D68').Value = 'ایشان که دارند مى‌آیند چه کسانى هستند؟'
$ws.Range('A69').Value = 'أشبل'
$ws.Range('B69').Value = 'مائي أشبل ٚ مرأ کۊکۊ چأکۊنيدي'
$ws.Range('C69').Value = 'تخم ماهى'
$ws.Range('D69').Value = 'با تخم ماهى کوکو درست مى‌کنند'
$ws.Range('A70').Value = 'أغۊز'
$ws.Range('B70').Value = 'أغۊز خۊري؟'
$ws.Range('C70').Value = 'گردو'
$ws.Range('D70').Value = 'گردو می‌خوری؟'
$ws.Range('B106').Value = 'أن ٚ باقي بمانسته    ٰ چي بۊکۊنم؟'
$ws.Range('B174').Value = 'أغۊز بگنسته مي سرأ'
$ws.Range('B189').Value = 'أشان خأئيد همه    ٰ بۊکۊشيد'
$ws.Range('B206').Value = 'أمي أغۊز دار هني پاچه'
$ws.Range('B209').Value = 'أمي مؤبل ٚ دسته    ٰ پت بزه'
$ws.Range('B298').Value = 'پارچهˈ همهˈ تۊرتۊره     ٰ کۊد'
$ws.Range('B316').Value = 'تي وسيله    ٰن کابينت ٚ جؤر نهأ'
$ws.Range('B330').Value = 'خۊ وسيله    ٰنأ جيجاماما کۊدن دره'
$ws.Range('B359').Value = 'پارچه     ٰ چارکۊنج واوين'
$ws.Range('B371').Value = 'چمچه    ٰ بيدئي؟'
$ws.Range('B372').Value = 'تۊ چنتا أغۊز دأري؟'
$ws.Range('B487').Value = 'تي أغۊزان ديمه بنن'
$ws.Range('B493').Value = 'ديرۊ شيمي کۊچه    ٰ دٚوارستم'
$ws.Range('B562').Value = 'تي وسيله    ٰنأ بنه ميز ٚ سر'
$ws.Range('B586').Value = 'سٚفره     ٰ بۊشؤستم'
$ws.Range('A611').Value = 'شلار دأئن'
$ws.Range('B665').Value = 'سٚفره     ٰ فلگان'
$ws.Range('B720').Value = 'بۊشۊ اۊ کلانه    ٰ بأور'
$ws.Range('B724').Value = 'بۊشۊ کله    ٰ وأگيران'
$ws.Range('B751').Value = 'جغله    ٰ ببر گاره دۊرۊن بۊخۊسان'
$ws.Range('B768').Value = 'گرزه خۊراک أغۊزه'
$ws.Range('B773').Value = 'گمج ٚ مئن أغۊز قاتؤق چأکۊن'
$ws.Range('B799').Value = 'اۊ لته    ٰ مرأ فأدن ميزأ پأکۊنم'
$ws.Range('B802').Value = 'لپه    ٰ تأود'
$ws.Range('B811').Value = 'مئوه    ٰنأ بۊشؤستي؟'
$ws.Range('B848').Value = 'أمين خۊ مچه    ٰ عمل بۊکۊده'
$ws.Range('B926').Value = 'ايپچه اۊشنتر همه    ٰ وارسي کۊديد'
$ws.Range('B944').Value = 'گاب خۊ مانده    ٰ واليشتن دره'
